$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price cells so they stay text (matching original inlineStr type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "61.419.27"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "3.393.38"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "403.23"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.680"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -6.99%  "
$ws.Range("D11").Value = "41.57"
$ws.Range("E11").Value = "  -4.71%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "8.33"
$ws.Range("E13").Value = "  -5.91%  "
$ws.Range("D14").Value = "19.66"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "3.437.61"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "11.62"
$ws.Range("E16").Value = "  +6.26%  "
$ws.Range("D17").Value = "61.494.52"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "0.0000139"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").Value = "3.14"
$ws.Range("E20").Value = "  -6.10%  "
$ws.Range("D21").Value = "82.82"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "310.27"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "4.82"
$ws.Range("E25").Value = "  +9.94%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "29.40"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  +5.09%  "
$ws.Range("D28").Value = "7.98"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").Value = "2.73"
$ws.Range("E29").Value = "  +5.05%  "
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "43.16"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "11.29"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").Value = "0.0479"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").Value = "51.28"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -6.19%  "
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("D40").Value = "0.318"
$ws.Range("E40").Value = "  +10.15%  "
$ws.Range("D41").Value = "138.90"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "1.95"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "3.93"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "16.48"
$ws.Range("E45").Value = "  -6.10%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "21.16"
$ws.Range("E47").Value = "  -4.63%  "
$ws.Range("D48").Value = "2.093.78"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  +19.26%  "
$ws.Range("E51").Value = "  +2.12%  "
